$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Jan 09 17:46:42 EST 2023"
$ws.Range("B3").Value = "Mon Jan 09 17:46:52 EST 2023"
$ws.Range("B4").Value = "Mon Jan 09 17:47:02 EST 2023"
$ws.Range("B5").Value = "Mon Jan 09 17:47:12 EST 2023"
$ws.Range("B6").Value = "Mon Jan 09 17:47:22 EST 2023"
$ws.Range("B7").Value = "Mon Jan 09 17:47:33 EST 2023"
